# This script applies the "harvard case classification" update.
# The *_old columns (Ada_old, Avey_old, Babylon_old, Buoy_old, K health_old,
# WebMD_old, doctor_MA_old, doctor_NJ_old, doctor_TH_old) are recalculated with
# the new classification, and the average_doctor / average_doctor_old header
# labels are swapped to reflect the new "current" vs "old" doctor baseline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the average_doctor / average_doctor_old column headers (row 1)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4: update recalculated *_old values (and doctor averages)
$ws.Range("E4").Value = 0.48
$ws.Range("F4").Value = 0.052
$ws.Range("G4").Value = 0.229
$ws.Range("N4").Value = 0.461
$ws.Range("O4").Value = 0.063
$ws.Range("P4").Value = 0.25
$ws.Range("Q4").Value = 0.044
$ws.Range("R4").Value = 0.031
$ws.Range("S4").Value = 0.176
$ws.Range("W4").Value = 0.369
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.322
$ws.Range("AI4").Value = 0.493
$ws.Range("AJ4").Value = 0.07000000000000001
$ws.Range("AK4").Value = 0.264
$ws.Range("AU4").Value = 0.254
$ws.Range("AV4").Value = 0.027
$ws.Range("AW4").Value = 0.165
$ws.Range("BA4").Value = 2.015
$ws.Range("BB4").Value = 0.146
$ws.Range("BC4").Value = 0.382
$ws.Range("BG4").Value = 0.71
$ws.Range("BH4").Value = 0.143
$ws.Range("BI4").Value = 0.379
$ws.Range("BM4").Value = 0.729
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.252
$ws.Range("BP4").Value = 0.672
$ws.Range("BQ4").Value = 0.762

# Row 5: update recalculated *_old values (and doctor averages)
$ws.Range("E5").Value = 0.6
$ws.Range("F5").Value = 0.048
$ws.Range("G5").Value = 0.22
$ws.Range("N5").Value = 0.715
$ws.Range("O5").Value = 0.078
$ws.Range("P5").Value = 0.28
$ws.Range("Q5").Value = 0.022
$ws.Range("R5").Value = 0.006
$ws.Range("S5").Value = 0.075
$ws.Range("W5").Value = 0.335
$ws.Range("X5").Value = 0.096
$ws.Range("Y5").Value = 0.31
$ws.Range("AI5").Value = 0.494
$ws.Range("AJ5").Value = 0.065
$ws.Range("AK5").Value = 0.255
$ws.Range("AU5").Value = 0.468
$ws.Range("AV5").Value = 0.08
$ws.Range("AW5").Value = 0.282
$ws.Range("BA5").Value = 1.251
$ws.Range("BB5").Value = 0.07000000000000001
$ws.Range("BC5").Value = 0.264
$ws.Range("BG5").Value = 0.357
$ws.Range("BH5").Value = 0.039
$ws.Range("BI5").Value = 0.198
$ws.Range("BM5").Value = 0.509
$ws.Range("BN5").Value = 0.039
$ws.Range("BO5").Value = 0.198
$ws.Range("BP5").Value = 0.417
$ws.Range("BQ5").Value = 0.442

# Row 6: update recalculated *_old values (and doctor averages)
$ws.Range("E6").Value = 0.533
$ws.Range("N6").Value = 0.5610000000000001
$ws.Range("Q6").Value = 0.029
$ws.Range("W6").Value = 0.351
$ws.Range("AI6").Value = 0.493
$ws.Range("AU6").Value = 0.329
$ws.Range("BA6").Value = 1.536
$ws.Range("BG6").Value = 0.475
$ws.Range("BM6").Value = 0.599
$ws.Range("BP6").Value = 0.512
$ws.Range("BQ6").Value = 0.5570000000000001

# Row 7: update recalculated *_old values (and doctor averages)
$ws.Range("E7").Value = 0.571
$ws.Range("N7").Value = 0.644
$ws.Range("Q7").Value = 0.024
$ws.Range("W7").Value = 0.341
$ws.Range("AI7").Value = 0.494
$ws.Range("AU7").Value = 0.401
$ws.Range("BA7").Value = 1.35
$ws.Range("BG7").Value = 0.396
$ws.Range("BM7").Value = 0.542
$ws.Range("BP7").Value = 0.45
$ws.Range("BQ7").Value = 0.482

# Row 8: update recalculated *_old values (and doctor averages)
$ws.Range("E8").Value = 0.701
$ws.Range("F8").Value = 0.065
$ws.Range("G8").Value = 0.255
$ws.Range("N8").Value = 0.786
$ws.Range("O8").Value = 0.068
$ws.Range("P8").Value = 0.261
$ws.Range("Q8").Value = 0.023
$ws.Range("S8").Value = 0.114
$ws.Range("W8").Value = 0.419
$ws.Range("X8").Value = 0.116
$ws.Range("Y8").Value = 0.341
$ws.Range("AI8").Value = 0.58
$ws.Range("AJ8").Value = 0.105
$ws.Range("AK8").Value = 0.325
$ws.Range("AU8").Value = 0.425
$ws.Range("AV8").Value = 0.08599999999999999
$ws.Range("AW8").Value = 0.293
$ws.Range("BA8").Value = 1.717
$ws.Range("BB8").Value = 0.107
$ws.Range("BC8").Value = 0.328
$ws.Range("BG8").Value = 0.53
$ws.Range("BH8").Value = 0.109
$ws.Range("BI8").Value = 0.33
$ws.Range("BM8").Value = 0.673
$ws.Range("BN8").Value = 0.059
$ws.Range("BO8").Value = 0.244
$ws.Range("BP8").Value = 0.572
$ws.Range("BQ8").Value = 0.617

# Row 9: update recalculated *_old values (and doctor averages)
$ws.Range("E9").Value = 0.657
$ws.Range("F9").Value = 0.225
$ws.Range("G9").Value = 0.475
$ws.Range("N9").Value = 0.6860000000000001
$ws.Range("O9").Value = 0.216
$ws.Range("P9").Value = 0.464
$ws.Range("W9").Value = 0.314
$ws.Range("X9").Value = 0.216
$ws.Range("Y9").Value = 0.464
$ws.Range("AI9").Value = 0.514
$ws.Range("AJ9").Value = 0.25
$ws.Range("AK9").Value = 0.5
$ws.Range("BA9").Value = 1.657
$ws.Range("BB9").Value = 0.248
$ws.Range("BC9").Value = 0.498
$ws.Range("BG9").Value = 0.543
$ws.Range("BH9").Value = 0.248
$ws.Range("BI9").Value = 0.498
$ws.Range("BM9").Value = 0.657
$ws.Range("BN9").Value = 0.225
$ws.Range("BO9").Value = 0.475
$ws.Range("BP9").Value = 0.552
$ws.Range("BQ9").Value = 0.604

# Row 10: update recalculated *_old values (and doctor averages)
$ws.Range("E10").Value = 0.8
$ws.Range("F10").Value = 0.16
$ws.Range("G10").Value = 0.4
$ws.Range("N10").Value = 0.914
$ws.Range("O10").Value = 0.078
$ws.Range("P10").Value = 0.28
$ws.Range("W10").Value = 0.543
$ws.Range("AI10").Value = 0.629
$ws.Range("AJ10").Value = 0.233
$ws.Range("AK10").Value = 0.483
$ws.Range("AU10").Value = 0.429
$ws.Range("AV10").Value = 0.245
$ws.Range("AW10").Value = 0.495
$ws.Range("BA10").Value = 2.144
$ws.Range("BB10").Value = 0.216
$ws.Range("BC10").Value = 0.464
$ws.Range("BG10").Value = 0.629
$ws.Range("BH10").Value = 0.233
$ws.Range("BI10").Value = 0.483
$ws.Range("BM10").Value = 0.829
$ws.Range("BN10").Value = 0.142
$ws.Range("BO10").Value = 0.377
$ws.Range("BP10").Value = 0.715
$ws.Range("BQ10").Value = 0.757

# Row 11: update recalculated *_old values (and doctor averages)
$ws.Range("E11").Value = 0.857
$ws.Range("F11").Value = 0.122
$ws.Range("G11").Value = 0.35
$ws.Range("N11").Value = 0.914
$ws.Range("O11").Value = 0.078
$ws.Range("P11").Value = 0.28
$ws.Range("W11").Value = 0.543
$ws.Range("AI11").Value = 0.714
$ws.Range("AJ11").Value = 0.204
$ws.Range("AK11").Value = 0.452
$ws.Range("AU11").Value = 0.571
$ws.Range("AV11").Value = 0.245
$ws.Range("AW11").Value = 0.495
$ws.Range("BA11").Value = 2.144
$ws.Range("BB11").Value = 0.216
$ws.Range("BC11").Value = 0.464
$ws.Range("BG11").Value = 0.629
$ws.Range("BH11").Value = 0.233
$ws.Range("BI11").Value = 0.483
$ws.Range("BM11").Value = 0.829
$ws.Range("BN11").Value = 0.142
$ws.Range("BO11").Value = 0.377
$ws.Range("BP11").Value = 0.715
$ws.Range("BQ11").Value = 0.764

# Row 12: update recalculated *_old values (and doctor averages)
$ws.Range("E12").Value = 1.467
$ws.Range("F12").Value = 0.982
$ws.Range("G12").Value = 0.991
$ws.Range("N12").Value = 1.312
$ws.Range("O12").Value = 0.34
$ws.Range("P12").Value = 0.583
$ws.Range("W12").Value = 1.526
$ws.Range("X12").Value = 0.46
$ws.Range("Y12").Value = 0.678
$ws.Range("AI12").Value = 1.6
$ws.Range("AJ12").Value = 1.44
$ws.Range("AK12").Value = 1.2
$ws.Range("AU12").Value = 2.773
$ws.Range("AV12").Value = 3.63
$ws.Range("AW12").Value = 1.905
$ws.Range("BA12").Value = 3.881
$ws.Range("BB12").Value = 0.498
$ws.Range("BC12").Value = 0.706
$ws.Range("BG12").Value = 1.182
$ws.Range("BH12").Value = 0.24
$ws.Range("BI12").Value = 0.49
$ws.Range("BM12").Value = 1.241
$ws.Range("BN12").Value = 0.252
$ws.Range("BO12").Value = 0.502
$ws.Range("BP12").Value = 1.294
$ws.Range("BQ12").Value = 1.267

# Row 13: update recalculated *_old values (and doctor averages)
$ws.Range("E13").Value = 1.419
$ws.Range("F13").Value = 0.297
$ws.Range("G13").Value = 0.545
$ws.Range("N13").Value = 1.77
$ws.Range("O13").Value = 0.482
$ws.Range("P13").Value = 0.694
$ws.Range("W13").Value = 1.001
$ws.Range("X13").Value = 0.2
$ws.Range("Y13").Value = 0.448
$ws.Range("AI13").Value = 1.154
$ws.Range("AJ13").Value = 0.303
$ws.Range("AK13").Value = 0.551
$ws.Range("AU13").Value = 2.014
$ws.Range("AV13").Value = 0.336
$ws.Range("AW13").Value = 0.579
$ws.Range("BA13").Value = 2.104
$ws.Range("BB13").Value = 0.266
$ws.Range("BC13").Value = 0.515
$ws.Range("BG13").Value = 0.526
$ws.Range("BH13").Value = 0.045
$ws.Range("BI13").Value = 0.211
$ws.Range("BM13").Value = 0.78
$ws.Range("BN13").Value = 0.134
$ws.Range("BO13").Value = 0.365
$ws.Range("BP13").Value = 0.701
$ws.Range("BQ13").Value = 0.637
